# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (new report date 2022-01-07 / serial 44568)
# at the top of the data block, pushing the existing rows down by 3
# (A1:T438 -> A1:T441).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 421..423 - Excel shifts rows 421:438 down to 424:441
# and copies formatting (incl. the date style on column D) from the row above,
# matching how the existing rows are styled.
$ws.Range("421:423").Insert()

# Row 421 - Especial
$ws.Range("A421").Value = 5
$ws.Range("B421").Value = "Macroferia Regional de Talca"
$ws.Range("C421").Value = "Maule"
$ws.Range("D421").Value = 44568
$ws.Range("E421").Value = 7
$ws.Range("F421").Value = "Fruta"
$ws.Range("G421").Value = 100101
$ws.Range("H421").Value = "Berries"
$ws.Range("I421").Value = 100112025
$ws.Range("J421").Value = "Frutilla"
$ws.Range("K421").Value = "Sin especificar"
$ws.Range("L421").Value = "Especial"
$ws.Range("M421").Value = 230
$ws.Range("N421").Value = 7000
$ws.Range("O421").Value = 7000
$ws.Range("P421").Value = 7000
$ws.Range("Q421").Value = "$/bandeja 7 kilos"
$ws.Range("R421").Value = "Provincia de Melipilla"
$ws.Range("S421").Value = 1000
$ws.Range("T421").Value = 7

# Row 422 - Especial
$ws.Range("A422").Value = 5
$ws.Range("B422").Value = "Macroferia Regional de Talca"
$ws.Range("C422").Value = "Maule"
$ws.Range("D422").Value = 44568
$ws.Range("E422").Value = 7
$ws.Range("F422").Value = "Fruta"
$ws.Range("G422").Value = 100101
$ws.Range("H422").Value = "Berries"
$ws.Range("I422").Value = 100112025
$ws.Range("J422").Value = "Frutilla"
$ws.Range("K422").Value = "Sin especificar"
$ws.Range("L422").Value = "Especial"
$ws.Range("M422").Value = 200
$ws.Range("N422").Value = 7000
$ws.Range("O422").Value = 7000
$ws.Range("P422").Value = 7000
$ws.Range("Q422").Value = "$/caja 7 kilos"
$ws.Range("R422").Value = "Región del Maule"
$ws.Range("S422").Value = 1000
$ws.Range("T422").Value = 7

# Row 423 - Segunda
$ws.Range("A423").Value = 5
$ws.Range("B423").Value = "Macroferia Regional de Talca"
$ws.Range("C423").Value = "Maule"
$ws.Range("D423").Value = 44568
$ws.Range("E423").Value = 7
$ws.Range("F423").Value = "Fruta"
$ws.Range("G423").Value = 100101
$ws.Range("H423").Value = "Berries"
$ws.Range("I423").Value = 100112025
$ws.Range("J423").Value = "Frutilla"
$ws.Range("K423").Value = "Sin especificar"
$ws.Range("L423").Value = "Segunda"
$ws.Range("M423").Value = 150
$ws.Range("N423").Value = 5000
$ws.Range("O423").Value = 5000
$ws.Range("P423").Value = 5000
$ws.Range("Q423").Value = "$/bandeja 7 kilos"
$ws.Range("R423").Value = "Provincia de Melipilla"
$ws.Range("S423").Value = 714
$ws.Range("T423").Value = 7

Write-Output "inserted rows 421-423, new dimension should be A1:T441"
